# ISYS3001.docx edit: add the "Version management" paragraph plus an
# APA in-text citation (Sommerville, Software Engineering) referencing a
# new bibliography source, replacing the trailing "…" placeholder
# paragraph and the empty paragraph that follows it.

$d = $word.ActiveDocument

# --- 1. Register the bibliography source behind the citation -------------
# (Best-effort, mirrors what Word's "Insert Citation" does under the
# hood; some hosts do not persist this collection.)
try {
    $sourceXml = '<b:Source xmlns:b="http://schemas.openxmlformats.org/officeDocument/2006/bibliography"><b:Tag>Ian</b:Tag><b:SourceType>BookSection</b:SourceType><b:Guid>{4814DF66-BC91-40F1-9939-467BA3CF1049}</b:Guid><b:Title>Software Engineering</b:Title><b:Author><b:Author><b:NameList><b:Person><b:Last>Sommerville</b:Last><b:First>Ian</b:First></b:Person></b:NameList></b:Author><b:BookAuthor><b:NameList><b:Person><b:Last>Sommerville</b:Last><b:First>Ian</b:First></b:Person></b:NameList></b:BookAuthor></b:Author><b:BookTitle>Software engineering</b:BookTitle><b:RefOrder>1</b:RefOrder></b:Source>'
    $null = $d.Bibliography.Sources.Add($sourceXml)
} catch {
}

# --- 2. Find the placeholder "…" paragraph and replace its text ----------
$find = $d.Content.Find
$find.ClearFormatting()
$found = $find.Execute([char]0x2026, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $ellipsisPara = $find.Parent.Paragraphs(1)
    $pRange = $ellipsisPara.Range
    $bodyRange = $d.Range($pRange.Start, $pRange.End - 1)

    $sentenceXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Version management is the process of keeping track of different versions of software components and the systems in which these components are used</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>'
    $null = $bodyRange.InsertXML($sentenceXml)
}

# --- 3. Replace the trailing empty paragraph with the citation sentence --
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastRange = $lastPara.Range

$citationXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:sdt><w:sdtPr><w:id w:val="-1387640971"/><w:citation/></w:sdtPr><w:sdtContent><w:r><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:instrText xml:space="preserve"> CITATION Ian \l 3081 </w:instrText></w:r><w:r><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:noProof/></w:rPr><w:t>(Sommerville)</w:t></w:r><w:r><w:fldChar w:fldCharType="end"/></w:r></w:sdtContent></w:sdt><w:r><w:t>.</w:t></w:r></w:p>'
$null = $lastRange.InsertXML($citationXml)
